{"js": "// The document uses a template placeholder:\n//   employer.name.full(middle='full')\n// which the commit renames to the simpler helper call:\n//   employer.name_full()\n//\n// This string appears three times in the document body (including once\n// inside a table cell together with a trailing \".rstrip('.')\"). Word's\n// body.search() finds matches across paragraphs *and* table cells, and\n// range.insertText(..., \"Replace\") swaps the matched text in place while\n// keeping the surrounding run formatting intact.\n\nconst body = context.document.body;\nconst oldText = \"employer.name.full(middle=\\u2019full\\u2019)\";\nconst newText = \"employer.name_full()\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document uses a template placeholder:\n#   employer.name.full(middle='full')\n# which the commit renames to the simpler helper call:\n#   employer.name_full()\n#\n# This string appears three times in the document body (including once\n# inside a table cell together with a trailing \".rstrip('.')\"). A single\n# Find/Replace over $d.Content with Replace:=wdReplaceAll (2) covers every\n# story, including table cells, in one pass.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"employer.name.full(middle=\u2019full\u2019)\"\n$find.Replacement.Text = \"employer.name_full()\"\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
